# This workbook is a flat "productLine" lookup table exported from a
# MySQL/classicmodels-style query (productLine, textDescription,
# htmlDescription). The edit re-shapes it into a clean two-column
# "ProductLine" / "ProductLineDescription" table:
#   - drop the always-empty htmlDescription column (C)
#   - rename/recase the header row to ProductLine / ProductLineDescription
#   - highlight the header row in yellow so it stands out
#   - auto-fit column A to the new (shorter) header/content
#   - leave the selection where the author ended up (cell D11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop column C (htmlDescription) - it only ever held empty strings.
$ws.Columns.Item(3).Delete()

# 2. Rework the header row in place (row 1 stays row 1; only the text and
#    its styling change - the 7 data rows below are untouched).
$ws.Range("A1").Value = "ProductLine"
$ws.Range("B1").Value = "ProductLineDescription"

# 3. Give the header row a yellow fill so it reads as a header.
$ws.Range("A1:B1").Interior.Color = 65535

# 4. Resize column A to fit the new (shorter) header/content.
$ws.Range("A:A").EntireColumn.AutoFit()

# 5. Match the final selection left behind in the saved file.
[void]$ws.Range("D11").Select()
